$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 4995
$ws.Range("J10").Value = 4995
$ws.Range("L10").Value = 4995
$ws.Range("N10").Value = -5581
# Row 98
$ws.Range("H98").Value = 1497.0358
$ws.Range("I98").Value = 1108.6957
$ws.Range("K98").Value = 1108.6957
$ws.Range("M98").Value = 389.3043
# Row 112
$ws.Range("H112").Value = 6895.077
$ws.Range("J112").Value = 1893.6
$ws.Range("L112").Value = 5680.799999999999
$ws.Range("N112").Value = -7896.799999999999
# Row 122
$ws.Range("H122").Value = 1497.0358
$ws.Range("I122").Value = 1108.6957
$ws.Range("K122").Value = 3326.0871
$ws.Range("M122").Value = -876.0870999999997
# Row 125
$ws.Range("H125").Value = 5681.8184
$ws.Range("I125").Value = 812.8
$ws.Range("J125").Value = 7113.8823
$ws.Range("K125").Value = 7315.2
$ws.Range("L125").Value = 64024.9407
$ws.Range("M125").Value = -4855.2
$ws.Range("N125").Value = -68944.94070000001
# Row 132
$ws.Range("H132").Value = 604.6923
$ws.Range("I132").Value = 594.1316
$ws.Range("J132").Value = 1006
$ws.Range("K132").Value = 1782.3948
$ws.Range("L132").Value = 3018
$ws.Range("M132").Value = 747.6052
$ws.Range("N132").Value = -8078
# Row 137
$ws.Range("H137").Value = 682988.8
$ws.Range("I137").Value = 2805.4546
$ws.Range("J137").Value = 1237212.2
$ws.Range("K137").Value = 8416.363799999999
$ws.Range("L137").Value = 3711636.6
$ws.Range("M137").Value = -5866.363799999999
$ws.Range("N137").Value = -3716736.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17116.598
$ws.Range("I32").Value = 18519.316
$ws.Range("J32").Value = 5093.2856
$ws.Range("K32").Value = 18519.316
$ws.Range("L32").Value = 5093.2856
$ws.Range("M32").Value = -18232.316
$ws.Range("N32").Value = -5667.2856
# Row 74
$ws.Range("H74").Value = 3811.0789
$ws.Range("I74").Value = 1142.2122
$ws.Range("J74").Value = 21425.6
$ws.Range("K74").Value = 1142.2122
$ws.Range("L74").Value = 21425.6
$ws.Range("M74").Value = -268.2121999999999
$ws.Range("N74").Value = -23173.6
# Row 77
$ws.Range("H77").Value = 3811.0789
$ws.Range("I77").Value = 1142.2122
$ws.Range("J77").Value = 21425.6
$ws.Range("K77").Value = 5711.061
$ws.Range("L77").Value = 107128
$ws.Range("M77").Value = -1343.061
$ws.Range("N77").Value = -115864
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2451.0908
$ws.Range("I107").Value = 2393.1538
$ws.Range("J107").Value = 2534.7778
$ws.Range("K107").Value = 2393.1538
$ws.Range("L107").Value = 2534.7778
$ws.Range("M107").Value = -473.1538
$ws.Range("N107").Value = -6374.7778
# Row 132
$ws.Range("H132").Value = 52999.5
$ws.Range("J132").Value = 52999.5
$ws.Range("L132").Value = 52999.5
$ws.Range("N132").Value = -63119.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 134
$ws.Range("H134").Value = 2493.8572
$ws.Range("I134").Value = 2248.7368
$ws.Range("J134").Value = 2784.9375
$ws.Range("K134").Value = 6746.2104
$ws.Range("L134").Value = 8354.8125
$ws.Range("M134").Value = -4211.2104
$ws.Range("N134").Value = -13424.8125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4761.4707
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4761.4707
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14284.4121
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -15906.4121
# Row 71
$ws.Range("H71").Value = 4761.4707
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4761.4707
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 42853.2363
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -50965.2363
# Row 106
$ws.Range("H106").Value = 5857.143
$ws.Range("J106").Value = 5857.143
$ws.Range("L106").Value = 17571.429
$ws.Range("N106").Value = -19463.429

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5995.1665
$ws.Range("I80").Value = 8463.125
$ws.Range("J80").Value = 4020.8
$ws.Range("K80").Value = 8463.125
$ws.Range("L80").Value = 4020.8
$ws.Range("M80").Value = -7465.125
$ws.Range("N80").Value = -6016.8
# Row 83
$ws.Range("H83").Value = 5995.1665
$ws.Range("I83").Value = 8463.125
$ws.Range("J83").Value = 4020.8
$ws.Range("K83").Value = 42315.625
$ws.Range("L83").Value = 20104
$ws.Range("M83").Value = -37323.625
$ws.Range("N83").Value = -30088
# Row 107
$ws.Range("H107").Value = 6179
$ws.Range("I107").Value = 12212.5
$ws.Range("J107").Value = 1007.4286
$ws.Range("K107").Value = 12212.5
$ws.Range("L107").Value = 1007.4286
$ws.Range("M107").Value = -10292.5
$ws.Range("N107").Value = -4847.4286
# Row 135
$ws.Range("H135").Value = 60312.5
$ws.Range("J135").Value = 60312.5
$ws.Range("L135").Value = 60312.5
$ws.Range("N135").Value = -70452.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4451.75
$ws.Range("I7").Value = 4451.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4451.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4339.75
$ws.Range("N7").ClearContents()
# Row 40
$ws.Range("H40").Value = 3313.2173
$ws.Range("I40").Value = 3237.5
$ws.Range("J40").Value = 3486.2856
$ws.Range("K40").Value = 3237.5
$ws.Range("L40").Value = 3486.2856
$ws.Range("M40").Value = -3101.5
$ws.Range("N40").Value = -3758.2856
# Row 126
$ws.Range("H126").Value = 4451.75
$ws.Range("I126").Value = 4451.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13355.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10885.25
$ws.Range("N126").ClearContents()
# Row 133
$ws.Range("H133").Value = 47522.375
$ws.Range("J133").Value = 47522.375
$ws.Range("L133").Value = 47522.375
$ws.Range("N133").Value = -52582.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 13000
$ws.Range("I15").Value = 10000
$ws.Range("K15").Value = 10000
$ws.Range("M15").Value = -9712
# Row 122
$ws.Range("H122").Value = 3233.3333
$ws.Range("I122").Value = 3350
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10050
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7600
$ws.Range("N122").Value = -13900
